$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '54.408.83'
$ws.Range("E2").Value = '  -8.14%  '

# Row 3
$ws.Range("D3").Value = '2.422.13'
$ws.Range("E3").Value = '  -14.71%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '465.34'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -7.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.31%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.491'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -7.17%  '

# Row 9
$ws.Range("D9").Value = '2.437.84'
$ws.Range("E9").Value = '  -14.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0948'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.72%  '

# Row 11
$ws.Range("E11").Value = '  -9.61%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.319'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -8.75%  '

# Row 13
$ws.Range("E13").Value = '  -4.18%  '

# Row 14
$ws.Range("D14").Value = '2.836.33'
$ws.Range("E14").Value = '  -15.02%  '

# Row 15
$ws.Range("D15").Value = '54.194.85'
$ws.Range("E15").Value = '  -8.66%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -9.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000131'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.46%  '

# Row 18
$ws.Range("D18").Value = '2.428.47'
$ws.Range("E18").Value = '  -14.51%  '

# Row 19
$ws.Range("E19").Value = '  -11.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '312.83'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -11.50%  '

# Row 21
$ws.Range("E21").Value = '  -15.29%  '

# Row 22
$ws.Range("E22").Value = '  +0.35%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.69'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.91%  '

# Row 24
$ws.Range("E24").Value = '  -13.47%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '56.68'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -10.25%  '

# Row 26
$ws.Range("E26").Value = '  +0.29%  '

# Row 27
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.383'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -10.76%  '

# Row 28
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.154'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -10.49%  '

# Row 29
$ws.Range("D29").Value = '2.504.10'
$ws.Range("E29").Value = '  -15.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.17'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.94%  '

# Row 31
$ws.Range("E31").Value = '  -0.45%  '

# Row 32
$ws.Range("D32").Value = '0.0₃0714'
$ws.Range("E32").Value = '  -12.31%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.59'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.83%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.70'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -7.17%  '

# Row 35
$ws.Range("E35").Value = '  -12.60%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.22%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.53'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -15.78%  '

# Row 38
$ws.Range("E38").Value = '  -7.91%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.802'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -11.29%  '

# Row 40
$ws.Range("E40").Value = '  -8.13%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.993'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.62%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.603'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.40%  '

# Row 43
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.30'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.12%  '

# Row 44
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0529'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.83%  '

# Row 45
$ws.Range("E45").Value = '  -9.30%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.12'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.24%  '

# Row 47
$ws.Range("D47").Value = '1.973.65'
$ws.Range("E47").Value = '  -11.11%  '

# Row 48
$ws.Range("E48").Value = '  -2.97%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0868'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.42%  '

# Row 50
$ws.Range("E50").Value = '  -5.79%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.59'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -15.21%  '

